# Update stats for 2025-08 (row 21) in Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B21").Value = 6235
$ws.Range("C21").Value = 987
$ws.Range("D21").Value = 5584668
$ws.Range("E21").Value = 895.696551724138
$ws.Range("F21").Value = 8.227738239888915
$ws.Range("G21").Value = 4.113924050632911
$ws.Range("H21").Value = 27.47554652255253
